$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "a0"
$ws.Range("B3").Value = "a1"
$ws.Range("B4").Value = "a2"
$ws.Range("B5").Value = "b0"
$ws.Range("B6").Value = "b1"
$ws.Range("B7").Value = "c0"
$ws.Range("B8").Value = "c1"
$ws.Range("B9").Value = "c2"
$ws.Range("B10").Value = "c3"
